# cryptos.xlsx refresh — GitHub Actions price-scrape sync
# Updates the "Price" / "Volume(1h)" columns for every coin row with freshly
# scraped figures, and fixes four rows where the scraper re-sorted two
# neighbouring coins of (near-)identical rank (Uniswap/ShibaInu,
# PancakeSwap/ImmutableX, InjectiveProtocol/Toncoin, EnergySwap/Stellar) so
# the Coin/Link columns for those rows are swapped back in step with their
# Price/Volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per update: Row = sheet row, Coin/Link = $null when unchanged,
# Price/Volume = $null when that particular cell is untouched by this sync.
$updates = @(
    [pscustomobject]@{ Row = 2; Coin = $null; Link = $null; Price = '43.635.48'; Volume = '  -5.74%  ' },
    [pscustomobject]@{ Row = 3; Coin = $null; Link = $null; Price = '2.587.95'; Volume = '  -0.36%  ' },
    [pscustomobject]@{ Row = 4; Coin = $null; Link = $null; Price = '0.999'; Volume = '  -0.03%  ' },
    [pscustomobject]@{ Row = 5; Coin = $null; Link = $null; Price = '300.79'; Volume = '  -2.32%  ' },
    [pscustomobject]@{ Row = 6; Coin = $null; Link = $null; Price = '96.36'; Volume = '  -3.65%  ' },
    [pscustomobject]@{ Row = 7; Coin = $null; Link = $null; Price = '0.577'; Volume = '  -4.39%  ' },
    [pscustomobject]@{ Row = 8; Coin = $null; Link = $null; Price = $null; Volume = '  +0.09%  ' },
    [pscustomobject]@{ Row = 9; Coin = $null; Link = $null; Price = '0.556'; Volume = '  -4.33%  ' },
    [pscustomobject]@{ Row = 10; Coin = $null; Link = $null; Price = '37.10'; Volume = '  -5.64%  ' },
    [pscustomobject]@{ Row = 11; Coin = $null; Link = $null; Price = '0.0814'; Volume = '  -3.73%  ' },
    [pscustomobject]@{ Row = 12; Coin = $null; Link = $null; Price = '7.82'; Volume = '  -4.35%  ' },
    [pscustomobject]@{ Row = 13; Coin = $null; Link = $null; Price = '2.975.92'; Volume = '  -0.63%  ' },
    [pscustomobject]@{ Row = 14; Coin = $null; Link = $null; Price = $null; Volume = '  +1.06%  ' },
    [pscustomobject]@{ Row = 15; Coin = $null; Link = $null; Price = '2.581.17'; Volume = '  -0.88%  ' },
    [pscustomobject]@{ Row = 16; Coin = $null; Link = $null; Price = '0.888'; Volume = '  -4.02%  ' },
    [pscustomobject]@{ Row = 17; Coin = $null; Link = $null; Price = '14.36'; Volume = '  -4.64%  ' },
    [pscustomobject]@{ Row = 18; Coin = $null; Link = $null; Price = '43.638.86'; Volume = '  -6.00%  ' },
    [pscustomobject]@{ Row = 19; Coin = 'Uniswap'; Link = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price = '6.65'; Volume = '  -1.58%  ' },
    [pscustomobject]@{ Row = 20; Coin = 'ShibaInu'; Link = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price = '0.0₃0976'; Volume = '  -3.65%  ' },
    [pscustomobject]@{ Row = 21; Coin = $null; Link = $null; Price = '12.35'; Volume = '  -5.15%  ' },
    [pscustomobject]@{ Row = 22; Coin = $null; Link = $null; Price = '73.15'; Volume = '  +2.27%  ' },
    [pscustomobject]@{ Row = 23; Coin = $null; Link = $null; Price = '265.83'; Volume = '  -3.31%  ' },
    [pscustomobject]@{ Row = 24; Coin = 'PancakeSwap'; Link = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Price = '2.94'; Volume = '  -3.60%  ' },
    [pscustomobject]@{ Row = 25; Coin = 'ImmutableX'; Link = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price = '2.21'; Volume = '  +1.73%  ' },
    [pscustomobject]@{ Row = 26; Coin = $null; Link = $null; Price = '29.38'; Volume = '  -0.66%  ' },
    [pscustomobject]@{ Row = 27; Coin = $null; Link = $null; Price = $null; Volume = '  -0.01%  ' },
    [pscustomobject]@{ Row = 28; Coin = $null; Link = $null; Price = '10.24'; Volume = '  -3.65%  ' },
    [pscustomobject]@{ Row = 29; Coin = 'InjectiveProtocol'; Link = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; Price = '37.73'; Volume = '  -3.77%  ' },
    [pscustomobject]@{ Row = 30; Coin = 'Toncoin'; Link = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price = '2.15'; Volume = '  -7.28%  ' },
    [pscustomobject]@{ Row = 31; Coin = $null; Link = $null; Price = '5.98'; Volume = '  -5.46%  ' },
    [pscustomobject]@{ Row = 32; Coin = $null; Link = $null; Price = '3.62'; Volume = '  +0.00%  ' },
    [pscustomobject]@{ Row = 33; Coin = $null; Link = $null; Price = '2.23'; Volume = '  +1.46%  ' },
    [pscustomobject]@{ Row = 34; Coin = $null; Link = $null; Price = '151.80'; Volume = '  +0.95%  ' },
    [pscustomobject]@{ Row = 35; Coin = $null; Link = $null; Price = '2.79'; Volume = '  -1.43%  ' },
    [pscustomobject]@{ Row = 36; Coin = $null; Link = $null; Price = '0.0809'; Volume = '  -3.93%  ' },
    [pscustomobject]@{ Row = 37; Coin = $null; Link = $null; Price = $null; Volume = '  -4.80%  ' },
    [pscustomobject]@{ Row = 38; Coin = 'EnergySwap'; Link = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price = '24.39'; Volume = '  +5.62%  ' },
    [pscustomobject]@{ Row = 39; Coin = 'Stellar'; Link = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price = '0.121'; Volume = '  -1.69%  ' },
    [pscustomobject]@{ Row = 40; Coin = $null; Link = $null; Price = '16.77'; Volume = '  +4.74%  ' },
    [pscustomobject]@{ Row = 41; Coin = $null; Link = $null; Price = '3.49'; Volume = '  -4.14%  ' },
    [pscustomobject]@{ Row = 42; Coin = $null; Link = $null; Price = '0.0314'; Volume = '  -5.11%  ' },
    [pscustomobject]@{ Row = 43; Coin = $null; Link = $null; Price = '3.85'; Volume = '  -5.71%  ' },
    [pscustomobject]@{ Row = 44; Coin = $null; Link = $null; Price = '2.069.18'; Volume = '  -4.24%  ' },
    [pscustomobject]@{ Row = 45; Coin = $null; Link = $null; Price = '0.996'; Volume = '  -0.08%  ' },
    [pscustomobject]@{ Row = 46; Coin = $null; Link = $null; Price = '88.07'; Volume = '  -6.01%  ' },
    [pscustomobject]@{ Row = 47; Coin = $null; Link = $null; Price = '9.20'; Volume = '  -3.68%  ' },
    [pscustomobject]@{ Row = 48; Coin = $null; Link = $null; Price = '1.61'; Volume = '  +4.19%  ' },
    [pscustomobject]@{ Row = 49; Coin = $null; Link = $null; Price = '2.832.12'; Volume = '  -0.56%  ' },
    [pscustomobject]@{ Row = 50; Coin = $null; Link = $null; Price = '105.89'; Volume = '  -3.38%  ' },
    [pscustomobject]@{ Row = 51; Coin = $null; Link = $null; Price = $null; Volume = '  -5.21%  ' }
)

# Matches a plain decimal number (no thousands separators) — the shape
# Excel's own type-inference would coerce to a Number on entry. Prices such
# as "43.635.48" use '.' as a thousands separator and are never touched by
# this, but plain prices like "0.999" or "300.79" are, so those cells get
# pinned to Text format first to keep them stored the same way the
# original scrape stored them.
$numericPattern = '^[+-]?(\d+(\.\d*)?|\.\d+)$'

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.Coin -ne $null) {
        $ws.Cells.Item($row, 2).Value = $u.Coin
    }
    if ($u.Link -ne $null) {
        $ws.Cells.Item($row, 3).Value = $u.Link
    }
    if ($u.Price -ne $null) {
        $priceCell = $ws.Cells.Item($row, 4)
        if ($u.Price -match $numericPattern) {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $u.Price
    }
    if ($u.Volume -ne $null) {
        $ws.Cells.Item($row, 5).Value = $u.Volume
    }
}
